# Git101.pptx - "Turn back time with Git" slide:
# fill in the (previously empty) body placeholder with the
# "git revert HEAD" guidance.

$p = $ppt.ActivePresentation

# Slide 7 is "Turn back time with Git".
$slide = $p.Slides.Item(7)

# Shape 2 is "Content Placeholder 2" - the body text box that only held
# an empty paragraph before this edit.
$shape = $slide.Shapes.Item("Content Placeholder 2")

$tr = $shape.TextFrame.TextRange
$tr.Text = ""

# Paragraph 1 (top level): "If you just ran a commit/push to your repo:"
$run = $tr.InsertAfter("If you just ran a commit/push to your repo:")

# Paragraph 2 (indented one level): "git" + " revert HEAD" as two runs,
# same as the source slide which types "git" then appends the command.
$run = $run.InsertAfter("`r")
$run = $run.InsertAfter("git")
$run = $run.InsertAfter(" revert HEAD")

# Paragraph 3 (indented one level): trailing empty paragraph.
$run = $run.InsertAfter("`r")

# TextRange.IndentLevel is 1-based (1 = top level), so level-1 (XML
# <a:pPr lvl="1"/>) is IndentLevel 2.
$tr.Paragraphs(2).IndentLevel = 2
$tr.Paragraphs(3).IndentLevel = 2
